$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.838.90"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.873.07"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "301.10"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "0.5327"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("D8").Value = "0.3758"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "0.07175"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").Value = "21.63"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").Value = "0.8863"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "0.08104"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").Value = "1.873.20"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "93.12"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "5.271"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "14.73"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "0.000008548"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").Value = "1.000"
$ws.Range("D20").Value = "26.917.27"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").Value = "4.975"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").Value = "10.68"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "6.390"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").Value = "147.13"
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").Value = "2.250"
$ws.Range("E25").Value = "  -2.92%  "
$ws.Range("D26").Value = "1.730"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").Value = "114.38"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "4.745"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").Value = "4.572"
$ws.Range("E30").Value = "  -6.73%  "
$ws.Range("D31").Value = "0.09126"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").Value = "0.7977"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.989"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.172"
$ws.Range("E35").Value = "  -4.18%  "
$ws.Range("D36").Value = "0.5936"
$ws.Range("E36").Value = "  +3.61%  "
$ws.Range("D37").Value = "2.614"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").Value = "3.155"
$ws.Range("E38").Value = "  -6.15%  "
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").Value = "6.651"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").Value = "8.896"
$ws.Range("D43").Value = "115.89"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").Value = "0.5047"
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("D45").Value = "0.1495"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "9.922"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").Value = "1.620"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "37.64"
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("D50").Value = "0.06031"
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").Value = "62.23"
$ws.Range("E51").Value = "  -2.79%  "

$ws.Range("D2:D51").ClearFormats()
